$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update target-cluster labels (shared-string reshuffle: index 24 now
# "Inflammatory-Mac", index 25 now "MuSCs"; the former "Resolving-Mac" text
# is gone)
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("D9").Value = "MuSCs"

# Updated TPM-derived numeric values
$ws.Range("G2").Value = 0.027123
$ws.Range("H2").Value = 0.081369
$ws.Range("I2").Value = 0.0960827240265261
$ws.Range("J2").Value = 0.09608272402652611
$ws.Range("M2").Value = 0.1579376666666667
$ws.Range("N2").Value = 0.473813
$ws.Range("O2").Value = 0.05467876644486869
$ws.Range("P2").Value = 0.07340983674118848
$ws.Range("Q2").Value = 0.004283743333
$ws.Range("R2").Value = 0.038553689997
$ws.Range("S2").Value = 0.005253684826433193
$ws.Range("T2").Value = 0.007053417084435949
$ws.Range("G3").Value = 0.027123
$ws.Range("H3").Value = 0.081369
$ws.Range("I3").Value = 0.0960827240265261
$ws.Range("J3").Value = 0.09608272402652611
$ws.Range("O3").Value = 0.1724539210166233
$ws.Range("P3").Value = 0.2315307204300726
$ws.Range("Q3").Value = 0.013510698621
$ws.Range("R3").Value = 0.121596287589
$ws.Range("S3").Value = 0.01656984250033254
$ws.Range("T3").Value = 0.02224610231474543
$ws.Range("G4").Value = 0.027123
$ws.Range("H4").Value = 0.081369
$ws.Range("I4").Value = 0.0960827240265261
$ws.Range("J4").Value = 0.09608272402652611
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02135966666666667
$ws.Range("N4").Value = 0.064079
$ws.Range("O4").Value = 0.00739481752299059
$ws.Range("P4").Value = 0.0099280284174107
$ws.Range("Q4").Value = 0.0005793382389999999
$ws.Range("R4").Value = 0.005214044151
$ws.Range("S4").Value = 0.0007105142112880242
$ws.Range("T4").Value = 0.0009539120145575811
$ws.Range("G5").Value = 0.027123
$ws.Range("H5").Value = 0.081369
$ws.Range("I5").Value = 0.0960827240265261
$ws.Range("J5").Value = 0.09608272402652611
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.21104
$ws.Range("N5").Value = 4.422079999999999
$ws.Range("O5").Value = 0.7654724950155174
$ws.Range("P5").Value = 0.6851314144113283
$ws.Range("Q5").Value = 0.05997003791999999
$ws.Range("R5").Value = 0.3598202275199999
$ws.Range("S5").Value = 0.07354868248847234
$ws.Range("T5").Value = 0.06582929261278715
$ws.Range("I6").Value = 0.9039172759734738
$ws.Range("J6").Value = 0.9039172759734738
$ws.Range("M6").Value = 0.1579376666666667
$ws.Range("N6").Value = 0.473813
$ws.Range("O6").Value = 0.05467876644486869
$ws.Range("P6").Value = 0.07340983674118848
$ws.Range("Q6").Value = 0.040300164715
$ws.Range("R6").Value = 0.3627014824349999
$ws.Range("S6").Value = 0.04942508161843549
$ws.Range("T6").Value = 0.06635641965675253
$ws.Range("I7").Value = 0.9039172759734738
$ws.Range("J7").Value = 0.9039172759734738
$ws.Range("O7").Value = 0.1724539210166233
$ws.Range("P7").Value = 0.2315307204300726
$ws.Range("S7").Value = 0.1558840785162907
$ws.Range("T7").Value = 0.2092846181153271
$ws.Range("I8").Value = 0.9039172759734738
$ws.Range("J8").Value = 0.9039172759734738
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02135966666666667
$ws.Range("N8").Value = 0.064079
$ws.Range("O8").Value = 0.00739481752299059
$ws.Range("P8").Value = 0.0099280284174107
$ws.Range("Q8").Value = 0.005450239344999999
$ws.Range("R8").Value = 0.049052154105
$ws.Range("S8").Value = 0.006684303311702565
$ws.Range("T8").Value = 0.008974116402853118
$ws.Range("I9").Value = 0.9039172759734738
$ws.Range("J9").Value = 0.9039172759734738
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.21104
$ws.Range("N9").Value = 4.422079999999999
$ws.Range("O9").Value = 0.7654724950155174
$ws.Range("P9").Value = 0.6851314144113283
$ws.Range("Q9").Value = 0.5641800215999999
$ws.Range("R9").Value = 3.385080129599999
$ws.Range("S9").Value = 0.691923812527045
$ws.Range("T9").Value = 0.619302121798541
